# Update "feriados" (holidays) worksheet: the Data column (A2:A10) used to
# hold real date values; it is now stored as plain text strings
# ("YYYY-MM-DD") instead, so the number format changes from a date format
# to Text ("@" / numFmtId 49).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

$dates = @(
    "2016-01-01",
    "2016-02-03",
    "2016-04-07",
    "2016-05-01",
    "2016-06-25",
    "2016-09-07",
    "2016-09-25",
    "2016-10-04",
    "2016-12-25"
)

# Re-format column A (rows 2-10) as Text, then write the ISO date strings.
$ws.Range("A2:A10").NumberFormat = "@"
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
}

# C2 loses its (redundant) explicit font/bold formatting, keeping only the
# centered alignment shared by the other description cells.
$ws.Range("C2").HorizontalAlignment = -4108

# Selection left on the whole of column A, matching the refreshed view.
$null = $ws.Columns.Item(1).Select()
